$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Regime Atual
$ws.Range("B2").Value = 0.6305335119698572
$ws.Range("C2").Value = 0.1351436301569832
$ws.Range("D2").Value = 0.5800723306279203
$ws.Range("E2").Value = 0.330484099728426
$ws.Range("F2").Value = 352.8774543083251

# Row 3 - Nova Proposta
$ws.Range("B3").Value = 0.6313367221965587
$ws.Range("C3").Value = 0.13443602830935
$ws.Range("D3").Value = 0.5808776223237405
$ws.Range("E3").Value = 0.3277835981368859
$ws.Range("F3").Value = 317.1246891006634
$ws.Range("G3").Value = -35.75276520766164

# Row 4 - Nova c/ Aliq. Maxima
$ws.Range("B4").Value = 0.6313367221965587
$ws.Range("C4").Value = 0.13443602830935
$ws.Range("D4").Value = 0.5808776223237405
$ws.Range("E4").Value = 0.3277835981368859
$ws.Range("F4").Value = 385.8038789611176
$ws.Range("G4").Value = 32.92642465279249
